$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> (D, J, K, L, M, P) target values after the edit.
# Row 8 is unchanged; row 1 (header) untouched.
$data = @{
    2  = @(44203, 30,  2000, 2000, 2000, 2000)
    3  = @(44452, 120, 2300, 2300, 2300, 2300)
    4  = @(44483, 50,  2200, 2200, 2200, 2200)
    5  = @(44487, 50,  2200, 2200, 2200, 2200)
    6  = @(44484, 40,  2200, 2200, 2200, 2200)
    7  = @(44476, 30,  2200, 2200, 2200, 2200)
    9  = @(44473, 140, 1600, 1600, 1600, 1600)
    10 = @(44497, 50,  2200, 2200, 2200, 2200)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 4).Value  = $vals[0]   # D - Fecha
    $ws.Cells.Item($row, 10).Value = $vals[1]   # J - Volumen
    $ws.Cells.Item($row, 11).Value = $vals[2]   # K - Precio minimo
    $ws.Cells.Item($row, 12).Value = $vals[3]   # L - Precio maximo
    $ws.Cells.Item($row, 13).Value = $vals[4]   # M - Precio promedio ponderado
    $ws.Cells.Item($row, 16).Value = $vals[5]   # P - Precio $/Kg
}
